$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 24-26 were blank placeholder rows (style matching rows above but no
# data). Clone the formatting from the last filled-in diary row (23) down
# onto the three new rows, then fill in the three new diary entries.
$ws.Range("A23:G23").Copy() | Out-Null
$ws.Range("A24:A26").PasteSpecial(-4122) | Out-Null

# --- Row 24: Feb 19 2020 ---
$ws.Range("A24").Value = 43880
$ws.Range("B24").Value = "2:30PM - 5:00 PM"
$ws.Range("C24").Value = "Chris Zhang, Nicolas Grantham, and Hyun Jay Yang"
$ws.Range("D24").Value = "Finish and resubmit assignment 2"
$ws.Range("E24").Value = "We rewrote the document using higher level abstractions and went our discovery process more throughly"
$ws.Range("F24").Value = "Structuring the document into paragraphs and making links to our diagrams makes our document more understandable"
$ws.Range("G24").Value = "Feeling good overall"

# --- Row 26: Feb 21 2020 ---
$ws.Range("A26").Value = 43882
$ws.Range("B26").Value = "2:00PM - 7:10 PM"
$ws.Range("C26").Value = "Chris Zhang, Nicolas Grantham, and Hyun Jay Yang"
$ws.Range("G26").Value = "Exhausted"
$ws.Range("E26").Value = "Explained the social context, identified interesting pull requests and issues, and explained the architecture of our project in a concise document"
$ws.Range("D26").Value = "Finish and deliver our project assignment #4"
$ws.Range("F26").Value = "Since we had already worked on the essential features, we already had a general understanding of our project, so it was easier to understand the architecture because we knew the exact routes we had to study"

# --- Row 25: Feb 20 2020 ---
$ws.Range("A25").Value = 43881
$ws.Range("B25").Value = "5:00PM - 7:00 PM"
$ws.Range("C25").Value = "N/A"
$ws.Range("D25").Value = "Learn new expert key practices, what is social context, and how does architecture can help understanding code"
$ws.Range("E25").Value = "Understood what social context is and how it might affect our decisions when choosing a project, and  how professionals use architecture as a comprehension tool"
$ws.Range("G25").Value = "Feeling ok"
$ws.Range("F25").Value = "It was good to hear our guest speakers give suggestions on how to introduce new members to the project and team, and how they guide them in the proper way to contribute"

# Row heights grow with the extra wrapped text (matches Excel's row autofit
# for the new multi-line entries).
$ws.Rows.Item(24).RowHeight = 63
$ws.Rows.Item(25).RowHeight = 94.5
$ws.Rows.Item(26).RowHeight = 110.25

# The author scrolled down and left the cursor on the first new diary row.
$ws.Range("A24").Select() | Out-Null
